$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Two new test cases were written for the "Create" row, bumping the
# "Total Test Cases" count in C4 from 8 to 10. The dependent SUM/ratio
# formulas in G4 and G6 recalculate automatically.
$ws.Range("C4").Value = 10

# Reflect the resulting active cell/selection in the saved view.
$ws.Range("D6").Select()
